# Insert a new data row before the current row 23. This pushes the
# existing rows 23-30 down to 24-31 (matching the diff, which shows every
# row from 23 through 30 shifting its data down by one row, and a brand
# new row 31 appearing with what used to be row 30's content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with its data.
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44529
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112022
$ws.Range("G23").Value = "Arveja Verde"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12500
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 500
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
